$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.2
$ws.Range("H3").Value = 6
$ws.Range("K3").Value = 2.6
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 2.08
$ws.Range("S3").Value = 1.3
$ws.Range("T3").Value = 3.25
$ws.Range("W3").Value = 6.5
$ws.Range("X3").Value = 5.5
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 12
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 51
$ws.Range("AK3").Value = 201
$ws.Range("AL3").Value = 101
$ws.Range("AP3").Value = 21
$ws.Range("AQ3").Value = 13
$ws.Range("AT3").Value = 3.25
$ws.Range("AV3").Value = 81
$ws.Range("AY3").Value = 51

# Row 4 updates
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.25
$ws.Range("K4").Value = 2.25
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.25
$ws.Range("P4").Value = 3.75
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("Z4").Value = 12
$ws.Range("AA4").Value = 13
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 7.5
$ws.Range("AG4").Value = 251
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 26
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 8.5
$ws.Range("AP4").Value = 19
$ws.Range("AQ4").Value = 26
$ws.Range("AS4").Value = 126
$ws.Range("AT4").Value = 3
$ws.Range("AZ4").Value = 101
